$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shifted values in column D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg) for rows 157-390
$ws.Cells.Item(157, 4).Value = 44915
$ws.Cells.Item(158, 4).Value = 44326
$ws.Cells.Item(158, 10).Value = 120
$ws.Cells.Item(159, 4).Value = 44284
$ws.Cells.Item(159, 10).Value = 180
$ws.Cells.Item(160, 4).Value = 44901
$ws.Cells.Item(160, 10).Value = 120
$ws.Cells.Item(161, 4).Value = 44795
$ws.Cells.Item(161, 10).Value = 130
$ws.Cells.Item(162, 4).Value = 44323
$ws.Cells.Item(162, 10).Value = 160
$ws.Cells.Item(163, 4).Value = 44203
$ws.Cells.Item(163, 10).Value = 120
$ws.Cells.Item(164, 4).Value = 44558
$ws.Cells.Item(165, 4).Value = 44235
$ws.Cells.Item(166, 4).Value = 44417
$ws.Cells.Item(167, 4).Value = 44504
$ws.Cells.Item(167, 10).Value = 160
$ws.Cells.Item(168, 4).Value = 44159
$ws.Cells.Item(169, 4).Value = 44210
$ws.Cells.Item(170, 4).Value = 44827
$ws.Cells.Item(170, 10).Value = 120
$ws.Cells.Item(171, 4).Value = 44694
$ws.Cells.Item(171, 10).Value = 160
$ws.Cells.Item(172, 4).Value = 44698
$ws.Cells.Item(172, 10).Value = 120
$ws.Cells.Item(173, 4).Value = 44267
$ws.Cells.Item(174, 4).Value = 44316
$ws.Cells.Item(174, 10).Value = 160
$ws.Cells.Item(175, 4).Value = 44546
$ws.Cells.Item(175, 10).Value = 180
$ws.Cells.Item(176, 4).Value = 44405
$ws.Cells.Item(176, 10).Value = 160
$ws.Cells.Item(177, 4).Value = 44160
$ws.Cells.Item(177, 10).Value = 230
$ws.Cells.Item(178, 4).Value = 44522
$ws.Cells.Item(179, 4).Value = 44484
$ws.Cells.Item(180, 4).Value = 44692
$ws.Cells.Item(181, 4).Value = 44329
$ws.Cells.Item(182, 4).Value = 44335
$ws.Cells.Item(182, 10).Value = 160
$ws.Cells.Item(183, 4).Value = 44530
$ws.Cells.Item(183, 10).Value = 120
$ws.Cells.Item(184, 4).Value = 44526
$ws.Cells.Item(184, 10).Value = 160
$ws.Cells.Item(185, 4).Value = 44302
$ws.Cells.Item(185, 10).Value = 130
$ws.Cells.Item(186, 4).Value = 44811
$ws.Cells.Item(186, 10).Value = 120
$ws.Cells.Item(187, 4).Value = 44586
$ws.Cells.Item(188, 4).Value = 44568
$ws.Cells.Item(189, 4).Value = 44277
$ws.Cells.Item(189, 10).Value = 160
$ws.Cells.Item(190, 4).Value = 44186
$ws.Cells.Item(191, 4).Value = 44195
$ws.Cells.Item(192, 4).Value = 44370
$ws.Cells.Item(192, 10).Value = 180
$ws.Cells.Item(193, 4).Value = 44334
$ws.Cells.Item(193, 10).Value = 190
$ws.Cells.Item(194, 4).Value = 44567
$ws.Cells.Item(194, 10).Value = 180
$ws.Cells.Item(195, 4).Value = 44274
$ws.Cells.Item(196, 4).Value = 44280
$ws.Cells.Item(196, 10).Value = 120
$ws.Cells.Item(197, 4).Value = 44343
$ws.Cells.Item(197, 10).Value = 180
$ws.Cells.Item(198, 4).Value = 44517
$ws.Cells.Item(198, 10).Value = 160
$ws.Cells.Item(199, 4).Value = 44214
$ws.Cells.Item(199, 10).Value = 110
$ws.Cells.Item(200, 4).Value = 44309
$ws.Cells.Item(201, 4).Value = 44308
$ws.Cells.Item(202, 4).Value = 44676
$ws.Cells.Item(203, 4).Value = 44649
$ws.Cells.Item(203, 10).Value = 160
$ws.Cells.Item(204, 4).Value = 44434
$ws.Cells.Item(204, 10).Value = 140
$ws.Cells.Item(205, 4).Value = 44487
$ws.Cells.Item(205, 10).Value = 160
$ws.Cells.Item(206, 4).Value = 44161
$ws.Cells.Item(207, 4).Value = 44355
$ws.Cells.Item(207, 10).Value = 180
$ws.Cells.Item(208, 4).Value = 44708
$ws.Cells.Item(209, 4).Value = 44384
$ws.Cells.Item(210, 4).Value = 44246
$ws.Cells.Item(210, 10).Value = 160
$ws.Cells.Item(211, 4).Value = 44749
$ws.Cells.Item(212, 4).Value = 44202
$ws.Cells.Item(212, 10).Value = 120
$ws.Cells.Item(213, 4).Value = 44488
$ws.Cells.Item(213, 10).Value = 150
$ws.Cells.Item(214, 4).Value = 44610
$ws.Cells.Item(215, 4).Value = 44572
$ws.Cells.Item(215, 10).Value = 160
$ws.Cells.Item(216, 4).Value = 44837
$ws.Cells.Item(216, 10).Value = 120
$ws.Cells.Item(217, 4).Value = 44278
$ws.Cells.Item(217, 10).Value = 130
$ws.Cells.Item(218, 4).Value = 44624
$ws.Cells.Item(218, 10).Value = 160
$ws.Cells.Item(219, 4).Value = 44327
$ws.Cells.Item(219, 10).Value = 190
$ws.Cells.Item(220, 4).Value = 44432
$ws.Cells.Item(220, 10).Value = 150
$ws.Cells.Item(221, 4).Value = 44454
$ws.Cells.Item(221, 10).Value = 160
$ws.Cells.Item(222, 4).Value = 44721
$ws.Cells.Item(222, 10).Value = 120
$ws.Cells.Item(223, 4).Value = 44252
$ws.Cells.Item(223, 10).Value = 160
$ws.Cells.Item(224, 4).Value = 44362
$ws.Cells.Item(224, 10).Value = 180
$ws.Cells.Item(225, 4).Value = 44475
$ws.Cells.Item(226, 4).Value = 44512
$ws.Cells.Item(226, 10).Value = 160
$ws.Cells.Item(227, 4).Value = 44614
$ws.Cells.Item(227, 10).Value = 230
$ws.Cells.Item(228, 4).Value = 44753
$ws.Cells.Item(228, 10).Value = 120
$ws.Cells.Item(229, 4).Value = 44518
$ws.Cells.Item(229, 10).Value = 160
$ws.Cells.Item(230, 4).Value = 44812
$ws.Cells.Item(230, 10).Value = 120
$ws.Cells.Item(231, 4).Value = 44209
$ws.Cells.Item(231, 10).Value = 160
$ws.Cells.Item(232, 4).Value = 44385
$ws.Cells.Item(233, 4).Value = 44403
$ws.Cells.Item(233, 10).Value = 180
$ws.Cells.Item(234, 4).Value = 44482
$ws.Cells.Item(235, 4).Value = 44897
$ws.Cells.Item(235, 10).Value = 160
$ws.Cells.Item(236, 4).Value = 44554
$ws.Cells.Item(236, 10).Value = 120
$ws.Cells.Item(237, 4).Value = 44603
$ws.Cells.Item(238, 4).Value = 44463
$ws.Cells.Item(238, 10).Value = 160
$ws.Cells.Item(239, 4).Value = 44792
$ws.Cells.Item(239, 10).Value = 120
$ws.Cells.Item(240, 4).Value = 44445
$ws.Cells.Item(241, 4).Value = 44174
$ws.Cells.Item(241, 10).Value = 180
$ws.Cells.Item(242, 4).Value = 44200
$ws.Cells.Item(242, 10).Value = 120
$ws.Cells.Item(243, 4).Value = 44602
$ws.Cells.Item(243, 10).Value = 130
$ws.Cells.Item(244, 4).Value = 44754
$ws.Cells.Item(244, 10).Value = 110
$ws.Cells.Item(245, 4).Value = 44503
$ws.Cells.Item(246, 4).Value = 44469
$ws.Cells.Item(246, 10).Value = 160
$ws.Cells.Item(247, 4).Value = 44711
$ws.Cells.Item(247, 10).Value = 260
$ws.Cells.Item(248, 4).Value = 44410
$ws.Cells.Item(248, 10).Value = 120
$ws.Cells.Item(249, 4).Value = 44509
$ws.Cells.Item(250, 4).Value = 44573
$ws.Cells.Item(250, 10).Value = 160
$ws.Cells.Item(251, 4).Value = 44560
$ws.Cells.Item(251, 10).Value = 180
$ws.Cells.Item(252, 4).Value = 44420
$ws.Cells.Item(253, 4).Value = 44636
$ws.Cells.Item(253, 10).Value = 160
$ws.Cells.Item(254, 4).Value = 44172
$ws.Cells.Item(254, 10).Value = 110
$ws.Cells.Item(255, 4).Value = 44301
$ws.Cells.Item(255, 10).Value = 130
$ws.Cells.Item(256, 4).Value = 44609
$ws.Cells.Item(256, 10).Value = 120
$ws.Cells.Item(257, 4).Value = 44579
$ws.Cells.Item(257, 10).Value = 160
$ws.Cells.Item(258, 4).Value = 44257
$ws.Cells.Item(258, 10).Value = 120
$ws.Cells.Item(259, 4).Value = 44476
$ws.Cells.Item(260, 4).Value = 44508
$ws.Cells.Item(261, 4).Value = 44635
$ws.Cells.Item(261, 10).Value = 160
$ws.Cells.Item(262, 4).Value = 44671
$ws.Cells.Item(262, 10).Value = 150
$ws.Cells.Item(263, 4).Value = 44336
$ws.Cells.Item(263, 10).Value = 160
$ws.Cells.Item(264, 4).Value = 44775
$ws.Cells.Item(265, 4).Value = 44802
$ws.Cells.Item(265, 10).Value = 120
$ws.Cells.Item(266, 4).Value = 44727
$ws.Cells.Item(266, 10).Value = 180
$ws.Cells.Item(267, 4).Value = 44585
$ws.Cells.Item(267, 10).Value = 160
$ws.Cells.Item(268, 4).Value = 44750
$ws.Cells.Item(269, 4).Value = 44655
$ws.Cells.Item(270, 4).Value = 44769
$ws.Cells.Item(271, 4).Value = 44705
$ws.Cells.Item(271, 10).Value = 120
$ws.Cells.Item(272, 4).Value = 44494
$ws.Cells.Item(272, 10).Value = 190
$ws.Cells.Item(273, 4).Value = 44232
$ws.Cells.Item(273, 10).Value = 120
$ws.Cells.Item(274, 4).Value = 44328
$ws.Cells.Item(275, 4).Value = 44466
$ws.Cells.Item(275, 10).Value = 160
$ws.Cells.Item(275, 11).Value = 1500
$ws.Cells.Item(275, 12).Value = 1500
$ws.Cells.Item(275, 13).Value = 1500
$ws.Cells.Item(275, 16).Value = 500
$ws.Cells.Item(276, 4).Value = 44179
$ws.Cells.Item(276, 10).Value = 48
$ws.Cells.Item(276, 11).Value = 2000
$ws.Cells.Item(276, 12).Value = 2000
$ws.Cells.Item(276, 13).Value = 2000
$ws.Cells.Item(276, 16).Value = 667
$ws.Cells.Item(277, 4).Value = 44448
$ws.Cells.Item(277, 10).Value = 160
$ws.Cells.Item(278, 4).Value = 44783
$ws.Cells.Item(278, 10).Value = 150
$ws.Cells.Item(279, 4).Value = 44595
$ws.Cells.Item(279, 10).Value = 230
$ws.Cells.Item(280, 4).Value = 44442
$ws.Cells.Item(280, 10).Value = 180
$ws.Cells.Item(281, 4).Value = 44657
$ws.Cells.Item(281, 10).Value = 120
$ws.Cells.Item(282, 4).Value = 44529
$ws.Cells.Item(283, 4).Value = 44498
$ws.Cells.Item(283, 10).Value = 160
$ws.Cells.Item(284, 4).Value = 44685
$ws.Cells.Item(284, 10).Value = 120
$ws.Cells.Item(285, 4).Value = 44418
$ws.Cells.Item(285, 10).Value = 150
$ws.Cells.Item(286, 4).Value = 44397
$ws.Cells.Item(286, 10).Value = 160
$ws.Cells.Item(287, 4).Value = 44599
$ws.Cells.Item(287, 10).Value = 168
$ws.Cells.Item(288, 4).Value = 44908
$ws.Cells.Item(288, 10).Value = 54
$ws.Cells.Item(289, 4).Value = 44600
$ws.Cells.Item(290, 4).Value = 44460
$ws.Cells.Item(290, 10).Value = 160
$ws.Cells.Item(291, 4).Value = 44319
$ws.Cells.Item(291, 10).Value = 190
$ws.Cells.Item(292, 4).Value = 44396
$ws.Cells.Item(292, 10).Value = 160
$ws.Cells.Item(293, 4).Value = 44608
$ws.Cells.Item(294, 4).Value = 44777
$ws.Cells.Item(294, 10).Value = 120
$ws.Cells.Item(295, 4).Value = 44263
$ws.Cells.Item(295, 10).Value = 180
$ws.Cells.Item(296, 4).Value = 44729
$ws.Cells.Item(296, 10).Value = 120
$ws.Cells.Item(297, 4).Value = 44342
$ws.Cells.Item(297, 10).Value = 260
$ws.Cells.Item(298, 4).Value = 44281
$ws.Cells.Item(299, 4).Value = 44532
$ws.Cells.Item(300, 4).Value = 44428
$ws.Cells.Item(301, 4).Value = 44785
$ws.Cells.Item(301, 10).Value = 160
$ws.Cells.Item(302, 4).Value = 44194
$ws.Cells.Item(302, 10).Value = 80
$ws.Cells.Item(302, 11).Value = 1500
$ws.Cells.Item(302, 12).Value = 1500
$ws.Cells.Item(302, 13).Value = 1500
$ws.Cells.Item(302, 16).Value = 500
$ws.Cells.Item(303, 4).Value = 44669
$ws.Cells.Item(303, 10).Value = 85
$ws.Cells.Item(303, 11).Value = 2000
$ws.Cells.Item(303, 12).Value = 2000
$ws.Cells.Item(303, 13).Value = 2000
$ws.Cells.Item(303, 16).Value = 667
$ws.Cells.Item(304, 4).Value = 44363
$ws.Cells.Item(305, 4).Value = 44322
$ws.Cells.Item(305, 10).Value = 130
$ws.Cells.Item(306, 4).Value = 44799
$ws.Cells.Item(307, 4).Value = 44344
$ws.Cells.Item(308, 4).Value = 44438
$ws.Cells.Item(309, 4).Value = 44426
$ws.Cells.Item(310, 4).Value = 44791
$ws.Cells.Item(311, 4).Value = 44511
$ws.Cells.Item(311, 10).Value = 160
$ws.Cells.Item(312, 4).Value = 44441
$ws.Cells.Item(312, 10).Value = 190
$ws.Cells.Item(313, 4).Value = 44638
$ws.Cells.Item(313, 10).Value = 180
$ws.Cells.Item(314, 4).Value = 44831
$ws.Cells.Item(314, 10).Value = 110
$ws.Cells.Item(315, 4).Value = 44270
$ws.Cells.Item(315, 10).Value = 120
$ws.Cells.Item(316, 4).Value = 44181
$ws.Cells.Item(316, 10).Value = 90
$ws.Cells.Item(317, 4).Value = 44406
$ws.Cells.Item(318, 4).Value = 44474
$ws.Cells.Item(318, 10).Value = 160
$ws.Cells.Item(319, 4).Value = 44651
$ws.Cells.Item(320, 4).Value = 44907
$ws.Cells.Item(321, 4).Value = 44663
$ws.Cells.Item(321, 10).Value = 120
$ws.Cells.Item(322, 4).Value = 44578
$ws.Cells.Item(322, 10).Value = 250
$ws.Cells.Item(323, 4).Value = 44490
$ws.Cells.Item(323, 10).Value = 160
$ws.Cells.Item(324, 4).Value = 44237
$ws.Cells.Item(325, 4).Value = 44566
$ws.Cells.Item(325, 10).Value = 130
$ws.Cells.Item(326, 4).Value = 44740
$ws.Cells.Item(327, 4).Value = 44895
$ws.Cells.Item(327, 10).Value = 120
$ws.Cells.Item(328, 4).Value = 44386
$ws.Cells.Item(329, 4).Value = 44425
$ws.Cells.Item(330, 4).Value = 44656
$ws.Cells.Item(330, 10).Value = 160
$ws.Cells.Item(331, 4).Value = 44211
$ws.Cells.Item(331, 10).Value = 120
$ws.Cells.Item(332, 4).Value = 44617
$ws.Cells.Item(332, 10).Value = 160
$ws.Cells.Item(333, 4).Value = 44452
$ws.Cells.Item(333, 10).Value = 190
$ws.Cells.Item(334, 4).Value = 44400
$ws.Cells.Item(334, 10).Value = 160
$ws.Cells.Item(334, 12).Value = 1500
$ws.Cells.Item(334, 13).Value = 1500
$ws.Cells.Item(334, 16).Value = 500
$ws.Cells.Item(335, 4).Value = 44559
$ws.Cells.Item(335, 10).Value = 172
$ws.Cells.Item(335, 12).Value = 2000
$ws.Cells.Item(335, 13).Value = 1747
$ws.Cells.Item(335, 16).Value = 582
$ws.Cells.Item(336, 4).Value = 44714
$ws.Cells.Item(336, 10).Value = 120
$ws.Cells.Item(337, 4).Value = 44293
$ws.Cells.Item(337, 10).Value = 160
$ws.Cells.Item(338, 4).Value = 44776
$ws.Cells.Item(338, 10).Value = 120
$ws.Cells.Item(339, 4).Value = 44491
$ws.Cells.Item(339, 10).Value = 160
$ws.Cells.Item(340, 4).Value = 44216
$ws.Cells.Item(340, 10).Value = 80
$ws.Cells.Item(341, 4).Value = 44264
$ws.Cells.Item(341, 10).Value = 120
$ws.Cells.Item(342, 4).Value = 44376
$ws.Cells.Item(342, 10).Value = 160
$ws.Cells.Item(343, 4).Value = 44305
$ws.Cells.Item(343, 10).Value = 180
$ws.Cells.Item(344, 4).Value = 44847
$ws.Cells.Item(344, 10).Value = 160
$ws.Cells.Item(345, 4).Value = 44839
$ws.Cells.Item(345, 10).Value = 120
$ws.Cells.Item(346, 4).Value = 44629
$ws.Cells.Item(346, 10).Value = 130
$ws.Cells.Item(347, 4).Value = 44592
$ws.Cells.Item(347, 10).Value = 160
$ws.Cells.Item(348, 4).Value = 44299
$ws.Cells.Item(348, 10).Value = 130
$ws.Cells.Item(349, 4).Value = 44382
$ws.Cells.Item(349, 10).Value = 160
$ws.Cells.Item(349, 11).Value = 1500
$ws.Cells.Item(349, 12).Value = 1500
$ws.Cells.Item(349, 13).Value = 1500
$ws.Cells.Item(349, 16).Value = 500
$ws.Cells.Item(350, 4).Value = 44756
$ws.Cells.Item(350, 10).Value = 78
$ws.Cells.Item(350, 11).Value = 1800
$ws.Cells.Item(350, 12).Value = 1800
$ws.Cells.Item(350, 13).Value = 1800
$ws.Cells.Item(350, 16).Value = 600
$ws.Cells.Item(351, 4).Value = 44725
$ws.Cells.Item(351, 10).Value = 180
$ws.Cells.Item(352, 4).Value = 44364
$ws.Cells.Item(353, 4).Value = 44818
$ws.Cells.Item(353, 10).Value = 160
$ws.Cells.Item(354, 4).Value = 44453
$ws.Cells.Item(354, 10).Value = 130
$ws.Cells.Item(355, 4).Value = 44516
$ws.Cells.Item(355, 10).Value = 150
$ws.Cells.Item(356, 4).Value = 44307
$ws.Cells.Item(356, 10).Value = 130
$ws.Cells.Item(357, 4).Value = 44162
$ws.Cells.Item(357, 10).Value = 160
$ws.Cells.Item(358, 4).Value = 44706
$ws.Cells.Item(359, 4).Value = 44883
$ws.Cells.Item(360, 4).Value = 44687
$ws.Cells.Item(360, 10).Value = 120
$ws.Cells.Item(361, 4).Value = 44468
$ws.Cells.Item(361, 10).Value = 180
$ws.Cells.Item(362, 4).Value = 44690
$ws.Cells.Item(363, 4).Value = 44533
$ws.Cells.Item(363, 10).Value = 160
$ws.Cells.Item(364, 4).Value = 44665
$ws.Cells.Item(364, 10).Value = 120
$ws.Cells.Item(365, 4).Value = 44477
$ws.Cells.Item(365, 10).Value = 160
$ws.Cells.Item(366, 4).Value = 44813
$ws.Cells.Item(367, 4).Value = 44217
$ws.Cells.Item(367, 10).Value = 120
$ws.Cells.Item(368, 4).Value = 44545
$ws.Cells.Item(368, 10).Value = 180
$ws.Cells.Item(369, 4).Value = 44819
$ws.Cells.Item(369, 10).Value = 160
$ws.Cells.Item(370, 4).Value = 44295
$ws.Cells.Item(370, 10).Value = 120
$ws.Cells.Item(370, 11).Value = 1500
$ws.Cells.Item(370, 12).Value = 1500
$ws.Cells.Item(370, 13).Value = 1500
$ws.Cells.Item(370, 16).Value = 500
$ws.Cells.Item(371, 4).Value = 44291
$ws.Cells.Item(371, 10).Value = 89
$ws.Cells.Item(371, 11).Value = 1800
$ws.Cells.Item(371, 12).Value = 1800
$ws.Cells.Item(371, 13).Value = 1800
$ws.Cells.Item(371, 16).Value = 600
$ws.Cells.Item(372, 4).Value = 44613
$ws.Cells.Item(372, 10).Value = 120
$ws.Cells.Item(373, 4).Value = 44630
$ws.Cells.Item(374, 4).Value = 44623
$ws.Cells.Item(374, 10).Value = 160
$ws.Cells.Item(375, 4).Value = 44679
$ws.Cells.Item(375, 10).Value = 180
$ws.Cells.Item(376, 4).Value = 44414
$ws.Cells.Item(376, 10).Value = 160
$ws.Cells.Item(377, 4).Value = 44741
$ws.Cells.Item(377, 10).Value = 120
$ws.Cells.Item(378, 4).Value = 44350
$ws.Cells.Item(379, 4).Value = 44447
$ws.Cells.Item(379, 10).Value = 160
$ws.Cells.Item(380, 4).Value = 44245
$ws.Cells.Item(380, 10).Value = 120
$ws.Cells.Item(381, 4).Value = 44580
$ws.Cells.Item(381, 10).Value = 160
$ws.Cells.Item(382, 4).Value = 44565
$ws.Cells.Item(382, 10).Value = 180
$ws.Cells.Item(383, 4).Value = 44589
$ws.Cells.Item(383, 10).Value = 150
$ws.Cells.Item(384, 4).Value = 44622
$ws.Cells.Item(384, 10).Value = 120
$ws.Cells.Item(385, 4).Value = 44571
$ws.Cells.Item(385, 10).Value = 190
$ws.Cells.Item(386, 4).Value = 44650
$ws.Cells.Item(386, 10).Value = 110
$ws.Cells.Item(387, 4).Value = 44806
$ws.Cells.Item(387, 10).Value = 50
$ws.Cells.Item(388, 4).Value = 44357
$ws.Cells.Item(388, 10).Value = 160
$ws.Cells.Item(389, 4).Value = 44911
$ws.Cells.Item(390, 4).Value = 44736
$ws.Cells.Item(390, 10).Value = 120

# Add new row 391 (duplicated/shifted from old row 390)
$ws.Cells.Item(391, 1).Value = 3
$ws.Cells.Item(391, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(391, 3).Value = 'Coquimbo'
$ws.Cells.Item(391, 4).Value = 44412
$ws.Cells.Item(391, 5).Value = 5
$ws.Cells.Item(391, 6).Value = 100112039
$ws.Cells.Item(391, 7).Value = 'Ciboulette'
$ws.Cells.Item(391, 8).Value = 'Sin especificar'
$ws.Cells.Item(391, 9).Value = 'Primera'
$ws.Cells.Item(391, 10).Value = 160
$ws.Cells.Item(391, 11).Value = 1500
$ws.Cells.Item(391, 12).Value = 1500
$ws.Cells.Item(391, 13).Value = 1500
$ws.Cells.Item(391, 14).Value = '$/docena de atados'
$ws.Cells.Item(391, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(391, 16).Value = 500
$ws.Cells.Item(391, 17).Value = 3
$ws.Cells.Item(391, 18).Value = 'Hortaliza'

# Fix number format for date column D391 to match other date cells
$ws.Cells.Item(391, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
